# feat: add 2022-Q3 data
#
# 1. Update the "总计" (totals) summary sheet with a new leading row for
#    2022-Q3 and shift the existing quarters down by one row.
# 2. Insert a brand-new "2022-Q3" worksheet (positioned right before the
#    existing "2022-Q2" sheet) carrying the per-fund holdings detail for
#    that quarter.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "总计" sheet: insert the 2022-Q3 row at the top of the data and
#    push 2022-Q2 / 2022-Q1 / 2021-Q4 down by one row each.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 4
$total.Range("D2").Value = 0.12

$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.02

$total.Range("B4").Value = "2022-Q1"
$total.Range("C4").Value = 5
$total.Range("D4").Value = 0.86

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q4"
$total.Range("C5").Value = 6
$total.Range("D5").Value = 0.67
# Row 5's "A" (index) cell needs the same look as A2:A4 (bold/centered/
# bordered style index already used on that column) - copy the format
# from the row above instead of hand-building a new style.
$total.Range("A4").Copy()
$total.Range("A5").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2. Brand-new "2022-Q3" sheet, inserted right before "2022-Q2".
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q3 = $wb.Worksheets.Add($q2)
$q3.Name = "2022-Q3"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $q3.Cells.Item(1, $i + 2).Value = $headers[$i]
}

$rows = @(
    @(0, "003300", "华夏圆和灵活配置混合A", "0.77", "75.31", "6.12", "0.0471", 7),
    @(1, "010746", "富安达长三角区域主题混合", "0.98", "88.86", "4.22", "0.0414", 9),
    @(2, "015068", "华夏圆和灵活配置混合C", "0.33", "75.31", "6.12", "0.0202", 7),
    @(3, "009789", "富安达科技创新混合", "0.46", "90.84", "3.44", "0.0158", 10)
)

for ($r = 0; $r -lt $rows.Count; $r++) {
    $row = $r + 2
    $data = $rows[$r]

    # 基金代码/基金规模/股票总仓位/仓位占比/持有市值 are stored as plain
    # text (not numbers) in every other quarter sheet - the fund code
    # in particular has significant leading zeros ("003300") that a
    # numeric cell would drop. Force a Text format before writing the
    # numeric-looking strings, otherwise they'd be auto-coerced into
    # numbers.
    $textCells = $q3.Range("B" + $row + ":G" + $row)
    $textCells.NumberFormat = "@"

    $q3.Cells.Item($row, 1).Value = $data[0]
    $q3.Cells.Item($row, 2).Value = $data[1]
    $q3.Cells.Item($row, 3).Value = $data[2]
    $q3.Cells.Item($row, 4).Value = $data[3]
    $q3.Cells.Item($row, 5).Value = $data[4]
    $q3.Cells.Item($row, 6).Value = $data[5]
    $q3.Cells.Item($row, 7).Value = $data[6]
    $q3.Cells.Item($row, 8).Value = $data[7]
}

# Re-apply the real formats by copying them in from an existing sheet
# that already has the right look, so the new sheet's style indices
# line up with the workbook's existing style table instead of minting
# fresh ones.
$q1 = $wb.Worksheets.Item("2022-Q1")

$q1.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)

$q1.Range("A2").Copy()
$q3.Range("A2:A5").PasteSpecial(-4122)

$q1.Range("B2").Copy()
$q3.Range("B2:G5").PasteSpecial(-4122)

$q1.Range("H2").Copy()
$q3.Range("H2:H5").PasteSpecial(-4122)

$q3.Range("A1").Select()
